$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Category names currently living in column A (rows 2-20), in order.
$names = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

# Insert a new column before column B. This shifts the old
# RawActivations / PercActivations / totalActivation columns (B, C, D) one
# step right (to C, D, E) and opens up a blank column B for the new
# "segments" label column.
$ws.Columns("B:B").Insert()

# Give the new B1 header the same bold/centered/bordered style already used
# by the other header cells (copy formatting from the now-shifted C1, which
# still carries it), then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B1").Value = "segments"

# The new segment-name cells (B2:B20) should be plain/unstyled like the rest
# of the data cells, not the bold/bordered style column A used to carry.
$ws.Range("B2:B20").Style = "Normal"

# Column A becomes a 0-based numeric index; column B gets the segment names
# that used to live in A.
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
}
